$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"

$ws.Range("G5").Select()
